$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds purely-numeric-looking text (e.g. "0.634", "15.10",
# "40.60", "0.0308"). Assigning such a string straight to .Value lets
# Excel auto-coerce it to a Double (dropping trailing zeros, or
# flipping to scientific notation for small values), which would
# corrupt the displayed text. Set-TextCell forces Text format first,
# writes the literal string, then restores the default "Normal"
# style so no stray NumberFormat is left on the cell (matches the
# original workbook, where these cells carry no explicit style).
function Set-TextCell {
    param($addr, $text)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextCell "D2" "42.056.47"
$ws.Range("E2").Value = "  -0.15%  "
Set-TextCell "D3" "2.230.51"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  -0.26%  "
Set-TextCell "D5" "252.31"
$ws.Range("E5").Value = "  +3.72%  "
Set-TextCell "D6" "0.634"
$ws.Range("E6").Value = "  +0.79%  "
Set-TextCell "D7" "71.58"
$ws.Range("E7").Value = "  +4.29%  "
$ws.Range("E8").Value = "  -0.07%  "
Set-TextCell "D9" "0.604"
$ws.Range("E9").Value = "  +9.87%  "
Set-TextCell "D10" "40.60"
$ws.Range("E10").Value = "  +13.38%  "
Set-TextCell "D11" "0.0974"
$ws.Range("E11").Value = "  -0.60%  "
Set-TextCell "D12" "58.28"
$ws.Range("E12").Value = "  -0.76%  "
Set-TextCell "D13" "7.29"
$ws.Range("E13").Value = "  +9.14%  "
Set-TextCell "D14" "0.104"
$ws.Range("E14").Value = "  -0.97%  "
Set-TextCell "D15" "2.552.74"
$ws.Range("E15").Value = "  -0.92%  "
Set-TextCell "D16" "15.10"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("E17").Value = "  +4.07%  "
Set-TextCell "D18" "2.218.54"
$ws.Range("E18").Value = "  -1.07%  "
Set-TextCell "D19" "41.886.47"
$ws.Range("E19").Value = "  -0.45%  "
Set-TextCell "D20" "0.0₃0970"
$ws.Range("E20").Value = "  +0.89%  "
Set-TextCell "D21" "6.31"
$ws.Range("E21").Value = "  +1.66%  "
Set-TextCell "D22" "72.93"
$ws.Range("E22").Value = "  -0.08%  "
Set-TextCell "D23" "235.93"
$ws.Range("E23").Value = "  +0.41%  "
Set-TextCell "D24" "2.11"
$ws.Range("E24").Value = "  +3.86%  "
Set-TextCell "D25" "4.12"
$ws.Range("E25").Value = "  +13.18%  "
Set-TextCell "D26" "11.96"
$ws.Range("E26").Value = "  +20.63%  "
$ws.Range("E27").Value = "  +0.04%  "
Set-TextCell "D28" "2.53"
$ws.Range("E28").Value = "  +2.83%  "
$ws.Range("E29").Value = "  -1.08%  "
Set-TextCell "D30" "170.83"
$ws.Range("E30").Value = "  -0.56%  "
Set-TextCell "D31" "20.94"
$ws.Range("E31").Value = "  +2.57%  "
Set-TextCell "D32" "0.123"
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D33" "5.62"
$ws.Range("E33").Value = "  +7.36%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D34" "0.125"
$ws.Range("E34").Value = "  -0.72%  "
Set-TextCell "D35" "0.0744"
$ws.Range("E35").Value = "  +4.33%  "
Set-TextCell "D36" "4.72"
$ws.Range("E36").Value = "  +0.65%  "
Set-TextCell "D37" "26.24"
$ws.Range("E37").Value = "  +17.34%  "
Set-TextCell "D38" "4.13"
$ws.Range("E38").Value = "  +9.95%  "
Set-TextCell "D39" "0.0308"
$ws.Range("E39").Value = "  +7.88%  "
Set-TextCell "D40" "2.30"
$ws.Range("E40").Value = "  +0.62%  "
Set-TextCell "D41" "5.94"
$ws.Range("E41").Value = "  +1.41%  "
Set-TextCell "D42" "12.56"
$ws.Range("E42").Value = "  +27.36%  "
Set-TextCell "D43" "66.22"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("E44").Value = "  +9.80%  "
Set-TextCell "D45" "4.85"
$ws.Range("E45").Value = "  -2.35%  "
Set-TextCell "D46" "8.81"
$ws.Range("E46").Value = "  -4.10%  "
Set-TextCell "D47" "0.103"
$ws.Range("E47").Value = "  +0.67%  "
Set-TextCell "D48" "4.68"
$ws.Range("E48").Value = "  +2.11%  "
$ws.Range("E49").Value = "  -0.09%  "
Set-TextCell "D50" "1.18"
$ws.Range("E50").Value = "  +7.64%  "
Set-TextCell "D51" "2.44"
$ws.Range("E51").Value = "  +5.85%  "
